$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as text (matches the source data's
# inline-string cells), even when the text looks like a plain number,
# without leaving a permanent "@" number-format on the cell.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$ws.Range("D2").Value = "61.019.54"
$ws.Range("E2").Value = "  -1.60%  "

$ws.Range("D3").Value = "3.411.93"
$ws.Range("E3").Value = "  -0.13%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.09%  "

Set-TextValue $ws.Range("D5") "574.68"
$ws.Range("E5").Value = "  -0.61%  "

Set-TextValue $ws.Range("D6") "138.13"
$ws.Range("E6").Value = "  -0.40%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.411.70"
$ws.Range("E8").Value = "  -0.16%  "

Set-TextValue $ws.Range("D9") "0.472"
$ws.Range("E9").Value = "  -1.21%  "

Set-TextValue $ws.Range("D10") "7.60"
$ws.Range("E10").Value = "  +1.40%  "

Set-TextValue $ws.Range("D11") "0.124"
$ws.Range("E11").Value = "  -2.70%  "

Set-TextValue $ws.Range("D12") "0.392"
$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("D13").Value = "3.989.85"
$ws.Range("E13").Value = "  -0.28%  "

$ws.Range("E14").Value = "  +1.06%  "

Set-TextValue $ws.Range("D15") "26.59"
$ws.Range("E15").Value = "  +4.06%  "

Set-TextValue $ws.Range("D16") "0.0000174"
$ws.Range("E16").Value = "  -3.34%  "

$ws.Range("D17").Value = "3.404.68"
$ws.Range("E17").Value = "  -0.41%  "

$ws.Range("D18").Value = "61.153.51"
$ws.Range("E18").Value = "  -1.43%  "

Set-TextValue $ws.Range("D19") "14.09"
$ws.Range("E19").Value = "  -0.68%  "

Set-TextValue $ws.Range("D20") "5.88"
$ws.Range("E20").Value = "  -0.50%  "

Set-TextValue $ws.Range("D21") "9.53"
$ws.Range("E21").Value = "  -0.22%  "

Set-TextValue $ws.Range("D22") "379.96"
$ws.Range("E22").Value = "  -3.01%  "

Set-TextValue $ws.Range("D23") "0.561"
$ws.Range("E23").Value = "  -2.08%  "

$ws.Range("D24").Value = "3.521.64"
$ws.Range("E24").Value = "  -0.89%  "

Set-TextValue $ws.Range("D25") "0.999"
$ws.Range("E25").Value = "  -0.22%  "

Set-TextValue $ws.Range("D26") "0.0000126"
$ws.Range("E26").Value = "  -1.67%  "

Set-TextValue $ws.Range("D27") "71.37"
$ws.Range("E27").Value = "  -0.23%  "

Set-TextValue $ws.Range("D28") "1.82"
$ws.Range("E28").Value = "  +14.26%  "

Set-TextValue $ws.Range("D29") "7.67"
$ws.Range("E29").Value = "  -0.01%  "

Set-TextValue $ws.Range("D30") "0.169"
$ws.Range("E30").Value = "  +4.86%  "

$ws.Range("E31").Value = "  +0.15%  "

Set-TextValue $ws.Range("D32") "8.24"
$ws.Range("E32").Value = "  -0.83%  "

Set-TextValue $ws.Range("D33") "2.17"
$ws.Range("E33").Value = "  -0.50%  "

$ws.Range("E34").Value = "  -0.02%  "

Set-TextValue $ws.Range("D35") "23.96"
$ws.Range("E35").Value = "  +1.47%  "

Set-TextValue $ws.Range("D36") "5.27"
$ws.Range("E36").Value = "  -4.19%  "

Set-TextValue $ws.Range("D37") "6.94"
$ws.Range("E37").Value = "  -1.13%  "

Set-TextValue $ws.Range("D38") "1.56"
$ws.Range("E38").Value = "  -1.25%  "

Set-TextValue $ws.Range("D39") "165.07"
$ws.Range("E39").Value = "  +1.53%  "

Set-TextValue $ws.Range("D40") "0.0765"
$ws.Range("E40").Value = "  -3.41%  "

$ws.Range("E41").Value = "  -0.04%  "

Set-TextValue $ws.Range("D42") "0.777"
$ws.Range("E42").Value = "  -1.98%  "

Set-TextValue $ws.Range("D43") "1.22"
$ws.Range("E43").Value = "  -1.26%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D44") "1.71"
$ws.Range("E44").Value = "  -3.21%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D45") "4.44"
$ws.Range("E45").Value = "  -0.77%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D46") "41.70"
$ws.Range("E46").Value = "  -0.03%  "

Set-TextValue $ws.Range("D47") "24.68"
$ws.Range("E47").Value = "  -2.48%  "

Set-TextValue $ws.Range("D48") "23.74"
$ws.Range("E48").Value = "  +2.44%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.429.01"
$ws.Range("E49").Value = "  +2.17%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D50") "6.84"
$ws.Range("E50").Value = "  -1.99%  "

Set-TextValue $ws.Range("D51") "2.48"
$ws.Range("E51").Value = "  +7.90%  "
